$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A14 with the new exception text (was "[Pintor Existente]", now
# "[Pintor já existe no sistema]") - same prefix/spacing, new wording.
$ws.Range("A14").Value = "Excepção 1               (passo 2)" + [char]10 + "[Pintor já existe no sistema]"

# The longer wrapped text needs a taller row (60 -> 90 points).
$ws.Range("A14").RowHeight = 90

# Author left the selection on D13 when they saved.
$ws.Range("D13").Select() | Out-Null
